# Add missing columns to tbl_Peixe in the "Atributos" sheet.
# Inserts 3 new rows right after the existing fk_id_localidade_peixe row
# (old row 45), pushing everything below it down by 3 rows, and fills
# them with the new "responsavel_projeto", "status" and
# "comprimento_padrão" column definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atributos")

# Insert three blank rows before row 45 (shifts existing rows 45+ down to 48+)
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# Row 45: responsavel_projeto
$ws.Cells.Item(45, 2).Value = "responsavel_projeto"
$ws.Cells.Item(45, 7).Value = "Pesquisador responsável pelo projeto em que a amostra foi coletada"
$ws.Cells.Item(45, 3).Value = "VARCHAR (255)"

# Row 46: status
$ws.Cells.Item(46, 2).Value = "status"
$ws.Cells.Item(46, 7).Value = "Status da amostra na coleção. Geralmente indica se a amostra acabou/ se esgotou"
$ws.Cells.Item(46, 3).Value = "VARCHAR (100)"

# Row 47: comprimento_padrão
$ws.Cells.Item(47, 2).Value = "comprimento_padrão"
$ws.Cells.Item(47, 7).Value = "Comprimento do espécime voucher em centimetros"
$ws.Cells.Item(47, 3).Value = "DECIMAL (5,2)"

# Update the sheet view to reflect where the author left the cursor
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 28
$ws.Range("A47").Select()
